# Auto commit: 2025-04-22 13:48:26
# Updates portfolio_statistics.xlsx: refresh Apple (column F) statistics with
# slightly revised floating point results, and append three new rows
# (CAGR, Max Drawdown, Calmar Ratio) to every *_Statistics worksheet.

$wb = $excel.ActiveWorkbook

$newRowLabels = @("CAGR", "Max Drawdown", "Calmar Ratio")

# ---- 1D_Statistics (sheet1) ----
$ws = $wb.Worksheets.Item("1D_Statistics")

# Refresh the APPLE (column F) figures with the recomputed values.
$ws.Cells.Item(2, 6).Value = 191.3809509277344
$ws.Cells.Item(3, 6).Value = 0.001754828062675263
$ws.Cells.Item(4, 6).Value = 0.01254982252972077
$ws.Cells.Item(5, 6).Value = 0.0001574980455274869
$ws.Cells.Item(6, 6).Value = -0.05249051910732176
$ws.Cells.Item(7, 6).Value = 1.470629158572999
$ws.Cells.Item(8, 6).Value = -0.01732569579896578
$ws.Cells.Item(9, 6).Value = -0.0253423453272444
$ws.Cells.Item(10, 6).Value = -0.02841676709174202
$ws.Cells.Item(11, 6).Value = -0.03844882503152471

# Append CAGR / Max Drawdown / Calmar Ratio rows, copying the number
# formatting from the row above so no new cell styles are introduced.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 1).Value = $newRowLabels[0]
$ws.Cells.Item(12, 2).Value = 0.5257815340198895
$ws.Cells.Item(12, 3).Value = 0.1165552785803567
$ws.Cells.Item(12, 4).Value = 0.104274499922256
$ws.Cells.Item(12, 5).Value = 0.5210181173153532
$ws.Cells.Item(12, 6).Value = 0.5431444910556542
$ws.Cells.Item(13, 1).Value = $newRowLabels[1]
$ws.Cells.Item(13, 2).Value = -0.1039095986491753
$ws.Cells.Item(13, 3).Value = -0.09780922394543157
$ws.Cells.Item(13, 4).Value = -0.1303597122302159
$ws.Cells.Item(13, 5).Value = -0.2793225480283115
$ws.Cells.Item(13, 6).Value = -0.1493236606462561
$ws.Cells.Item(14, 1).Value = $newRowLabels[2]
$ws.Cells.Item(14, 2).Value = 5.059990038023909
$ws.Cells.Item(14, 3).Value = 1.191659374021654
$ws.Cells.Item(14, 4).Value = 0.799898205805385
$ws.Cells.Item(14, 5).Value = 1.865292010949807
$ws.Cells.Item(14, 6).Value = 3.637363889319251

# ---- 3D_Statistics (sheet2) ----
$ws = $wb.Worksheets.Item("3D_Statistics")

# Refresh the APPLE (column F) figures with the recomputed values.
$ws.Cells.Item(2, 6).Value = 191.3809509277344
$ws.Cells.Item(3, 6).Value = 0.00526448418802579
$ws.Cells.Item(4, 6).Value = 0.02173693024744894
$ws.Cells.Item(5, 6).Value = 0.0004724941365824608
$ws.Cells.Item(6, 6).Value = -0.05249051910732176
$ws.Cells.Item(7, 6).Value = 1.470629158572999
$ws.Cells.Item(8, 6).Value = -0.03000898540029139
$ws.Cells.Item(9, 6).Value = -0.04389422968974303
$ws.Cells.Item(10, 6).Value = -0.04921928438974845
$ws.Cells.Item(11, 6).Value = -0.06659531844592684

# Append CAGR / Max Drawdown / Calmar Ratio rows, copying the number
# formatting from the row above so no new cell styles are introduced.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 1).Value = $newRowLabels[0]
$ws.Cells.Item(12, 2).Value = 0.5257815340198895
$ws.Cells.Item(12, 3).Value = 0.1165552785803567
$ws.Cells.Item(12, 4).Value = 0.104274499922256
$ws.Cells.Item(12, 5).Value = 0.5210181173153532
$ws.Cells.Item(12, 6).Value = 0.5431444910556542
$ws.Cells.Item(13, 1).Value = $newRowLabels[1]
$ws.Cells.Item(13, 2).Value = -0.1039095986491753
$ws.Cells.Item(13, 3).Value = -0.09780922394543157
$ws.Cells.Item(13, 4).Value = -0.1303597122302159
$ws.Cells.Item(13, 5).Value = -0.2793225480283115
$ws.Cells.Item(13, 6).Value = -0.1493236606462561
$ws.Cells.Item(14, 1).Value = $newRowLabels[2]
$ws.Cells.Item(14, 2).Value = 5.059990038023909
$ws.Cells.Item(14, 3).Value = 1.191659374021654
$ws.Cells.Item(14, 4).Value = 0.799898205805385
$ws.Cells.Item(14, 5).Value = 1.865292010949807
$ws.Cells.Item(14, 6).Value = 3.637363889319251

# ---- 5D_Statistics (sheet3) ----
$ws = $wb.Worksheets.Item("5D_Statistics")

# Refresh the APPLE (column F) figures with the recomputed values.
$ws.Cells.Item(2, 6).Value = 191.3809509277344
$ws.Cells.Item(3, 6).Value = 0.008774140313376318
$ws.Cells.Item(4, 6).Value = 0.02806225628201401
$ws.Cells.Item(5, 6).Value = 0.0007874902276374348
$ws.Cells.Item(6, 6).Value = -0.05249051910732176
$ws.Cells.Item(7, 6).Value = 1.470629158572999
$ws.Cells.Item(8, 6).Value = -0.03874143356397003
$ws.Cells.Item(9, 6).Value = -0.05666720686099264
$ws.Cells.Item(10, 6).Value = -0.06354182291791416
$ws.Cells.Item(11, 6).Value = -0.08597418642548475

# Append CAGR / Max Drawdown / Calmar Ratio rows, copying the number
# formatting from the row above so no new cell styles are introduced.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 1).Value = $newRowLabels[0]
$ws.Cells.Item(12, 2).Value = 0.5257815340198895
$ws.Cells.Item(12, 3).Value = 0.1165552785803567
$ws.Cells.Item(12, 4).Value = 0.104274499922256
$ws.Cells.Item(12, 5).Value = 0.5210181173153532
$ws.Cells.Item(12, 6).Value = 0.5431444910556542
$ws.Cells.Item(13, 1).Value = $newRowLabels[1]
$ws.Cells.Item(13, 2).Value = -0.1039095986491753
$ws.Cells.Item(13, 3).Value = -0.09780922394543157
$ws.Cells.Item(13, 4).Value = -0.1303597122302159
$ws.Cells.Item(13, 5).Value = -0.2793225480283115
$ws.Cells.Item(13, 6).Value = -0.1493236606462561
$ws.Cells.Item(14, 1).Value = $newRowLabels[2]
$ws.Cells.Item(14, 2).Value = 5.059990038023909
$ws.Cells.Item(14, 3).Value = 1.191659374021654
$ws.Cells.Item(14, 4).Value = 0.799898205805385
$ws.Cells.Item(14, 5).Value = 1.865292010949807
$ws.Cells.Item(14, 6).Value = 3.637363889319251

# ---- 10D_Statistics (sheet4) ----
$ws = $wb.Worksheets.Item("10D_Statistics")

# Refresh the APPLE (column F) figures with the recomputed values.
$ws.Cells.Item(2, 6).Value = 191.3809509277344
$ws.Cells.Item(3, 6).Value = 0.01754828062675264
$ws.Cells.Item(4, 6).Value = 0.0396860234248138
$ws.Cells.Item(5, 6).Value = 0.00157498045527487
$ws.Cells.Item(6, 6).Value = -0.05249051910732176
$ws.Cells.Item(7, 6).Value = 1.470629158572999
$ws.Cells.Item(8, 6).Value = -0.05478866077194264
$ws.Cells.Item(9, 6).Value = -0.08013953248461748
$ws.Cells.Item(10, 6).Value = -0.08986170774842375
$ws.Cells.Item(11, 6).Value = -0.1215858604569134

# Append CAGR / Max Drawdown / Calmar Ratio rows, copying the number
# formatting from the row above so no new cell styles are introduced.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 1).Value = $newRowLabels[0]
$ws.Cells.Item(12, 2).Value = 0.5257815340198895
$ws.Cells.Item(12, 3).Value = 0.1165552785803567
$ws.Cells.Item(12, 4).Value = 0.104274499922256
$ws.Cells.Item(12, 5).Value = 0.5210181173153532
$ws.Cells.Item(12, 6).Value = 0.5431444910556542
$ws.Cells.Item(13, 1).Value = $newRowLabels[1]
$ws.Cells.Item(13, 2).Value = -0.1039095986491753
$ws.Cells.Item(13, 3).Value = -0.09780922394543157
$ws.Cells.Item(13, 4).Value = -0.1303597122302159
$ws.Cells.Item(13, 5).Value = -0.2793225480283115
$ws.Cells.Item(13, 6).Value = -0.1493236606462561
$ws.Cells.Item(14, 1).Value = $newRowLabels[2]
$ws.Cells.Item(14, 2).Value = 5.059990038023909
$ws.Cells.Item(14, 3).Value = 1.191659374021654
$ws.Cells.Item(14, 4).Value = 0.799898205805385
$ws.Cells.Item(14, 5).Value = 1.865292010949807
$ws.Cells.Item(14, 6).Value = 3.637363889319251

# ---- Annual_Statistics (sheet5) ----
$ws = $wb.Worksheets.Item("Annual_Statistics")

# Refresh the APPLE (column F) figures with the recomputed values.
$ws.Cells.Item(2, 6).Value = 191.3809509277344
$ws.Cells.Item(3, 6).Value = 0.4387070156688158
$ws.Cells.Item(4, 6).Value = 0.198430117124069
$ws.Cells.Item(5, 6).Value = 0.03937451138187174
$ws.Cells.Item(6, 6).Value = -0.05249051910732176
$ws.Cells.Item(7, 6).Value = 1.470629158572999
$ws.Cells.Item(8, 6).Value = -0.2739433038597132
$ws.Cells.Item(9, 6).Value = -0.4006976624230874
$ws.Cells.Item(10, 6).Value = -0.4493085387421187
$ws.Cells.Item(11, 6).Value = -0.6079293022845669

# Append CAGR / Max Drawdown / Calmar Ratio rows, copying the number
# formatting from the row above so no new cell styles are introduced.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 1).Value = $newRowLabels[0]
$ws.Cells.Item(12, 2).Value = 0.5257815340198895
$ws.Cells.Item(12, 3).Value = 0.1165552785803567
$ws.Cells.Item(12, 4).Value = 0.104274499922256
$ws.Cells.Item(12, 5).Value = 0.5210181173153532
$ws.Cells.Item(12, 6).Value = 0.5431444910556542
$ws.Cells.Item(13, 1).Value = $newRowLabels[1]
$ws.Cells.Item(13, 2).Value = -0.1039095986491753
$ws.Cells.Item(13, 3).Value = -0.09780922394543157
$ws.Cells.Item(13, 4).Value = -0.1303597122302159
$ws.Cells.Item(13, 5).Value = -0.2793225480283115
$ws.Cells.Item(13, 6).Value = -0.1493236606462561
$ws.Cells.Item(14, 1).Value = $newRowLabels[2]
$ws.Cells.Item(14, 2).Value = 5.059990038023909
$ws.Cells.Item(14, 3).Value = 1.191659374021654
$ws.Cells.Item(14, 4).Value = 0.799898205805385
$ws.Cells.Item(14, 5).Value = 1.865292010949807
$ws.Cells.Item(14, 6).Value = 3.637363889319251
